# Se agrega modificacion para realizar deslogue y obtener los datos de prueba desde excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: test user "names" (PEPITO1..PEPITO5)
$ws.Range("A2").Value = "PEPITO1"
$ws.Range("A3").Value = "PEPITO2"
$ws.Range("A4").Value = "PEPITO3"
$ws.Range("A5").Value = "PEPITO4"
$ws.Range("A6").Value = "PEPITO5"

# Column B: test user "surnames" - PEREZ2 entered first for rows 3-6, then PEREZ1 for row 2
$ws.Range("B3").Value = "PEREZ2"
$ws.Range("B4").Value = "PEREZ2"
$ws.Range("B5").Value = "PEREZ2"
$ws.Range("B6").Value = "PEREZ2"
$ws.Range("B2").Value = "PEREZ1"

# Column C: reuse existing CARDONA value for all new rows
$ws.Range("C2").Value = "CARDONA"
$ws.Range("C3").Value = "CARDONA"
$ws.Range("C4").Value = "CARDONA"
$ws.Range("C5").Value = "CARDONA"
$ws.Range("C6").Value = "CARDONA"

# Leave the active cell/selection on B6, matching the saved worksheet state
[void]$ws.Range("B6").Select()
